# Re-run of the mtc2 experiment: refresh the neural-net sweep results
# (columns F:I, rows 2-16) with the latest training run's numbers, and
# drop the manual yellow/bordered highlight that had been left on the
# "[10, 10]" rows (13 and 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cost_trains (F), cost_vals (G), r2_trains (H), r2_vals (I) values
# for each data row, keyed by row number.
$newValues = @{
   2  = @(449.3388061523438,  20.26535224914551,  0.904532302286987,  0.780306371887003)
   3  = @(88.12675476074219,  23.34745216369629,  0.8687802853786691, 0.7468937940971623)
   4  = @(12.57365608215332,  18.68239212036133,  0.9048713433512923, 0.7974670287218527)
   5  = @(614.7761840820312,  14.07400989532471,  0.9326290046131448, 0.8474257942144309)
   6  = @(51.40799713134766,  14.06118488311768,  0.9291630328925761, 0.8475648216855438)
   7  = @(10.33803749084473,  17.25510787963867,  0.9217864186505018, 0.8129399897754346)
   8  = @(809.239013671875,   14.60351753234863,  0.9366141569452692, 0.8416854827648652)
   9  = @(76.40055084228516,  18.61332130432129,  0.9364223198641348, 0.7982157911542687)
   10 = @(8.794486999511719,  15.23386383056641,  0.9334647539550639, 0.834851986670313)
   11 = @(594.0044555664062,  12.20994853973389,  0.9647849320646769, 0.8676337874338909)
   12 = @(90.15464782714844,  13.22945976257324,  0.9476766135202127, 0.8565814078683964)
   13 = @(4.734439849853516,  10.21199321746826,  0.9641830124049209, 0.8892933200808956)
   14 = @(597.9210815429688,  11.91832733154297,  0.977185067894512,  0.8707952122350648)
   15 = @(96.18222045898438,  11.11820888519287,  0.9761843843674111, 0.8794691775277413)
   16 = @(3.886425971984863,  12.26212024688721,  0.9705983846100112, 0.8670681878577388)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    $ws.Cells.Item($r, 6).Value = $vals[0]   # F - cost_trains
    $ws.Cells.Item($r, 7).Value = $vals[1]   # G - cost_vals
    $ws.Cells.Item($r, 8).Value = $vals[2]   # H - r2_trains
    $ws.Cells.Item($r, 9).Value = $vals[3]   # I - r2_vals
}

# Rows 13 and 15 ("[10, 10]" architecture) previously had a manual
# yellow-fill + thin-border highlight applied (style index 2). Remove
# that highlight so the rows fall back to the sheet's default styling.
$ws.Range("A13:I13").ClearFormats()
$ws.Range("A15:I15").ClearFormats()

# Return the selection to the top-left cell (the highlighted K12
# selection left over from the previous save is no longer relevant).
[void]$ws.Range("A1").Select()

Write-Host "Applied mtc2 data refresh + highlight removal"
